$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 6 (existing rows 6-10 shift down to 7-11,
# carrying their formatting, values and row heights with them).
$ws.Rows.Item(6).Insert()

# Bump the "order" column (B) for every presentation that moved down one slot.
$ws.Range("B7").Value = 2
$ws.Range("B8").Value = 3
$ws.Range("B9").Value = 4
$ws.Range("B10").Value = 5
$ws.Range("B11").Value = 6

# Fill in the new CBIRT webinar entry in row 6.
$ws.Range("A6").Value = "conference"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = "Wright, J."
$ws.Range("D6").Value = "2022, October"
$ws.Range("E6").Value = "Translating Executive Functioning Challenges into Treatment Implementation to Support Return-to-Learn in the mTBI Population"
$ws.Range("P6").Value = "Center on Brain Injury Research and Training"
$ws.Range("Q6").Value = $false
$ws.Range("R6").Value = "Webinar "
$ws.Range("S6").Value = $false

# The row-above's format carried an unstyled placeholder cell in column J and
# no leftover cells in L/M - match that for row 6.
$ws.Range("J6").Style = "Normal"
$ws.Range("L6").Clear()
$ws.Range("M6").Clear()

# The new row's wrapped text renders at a height of 119 (matching the other
# 4-line entries in this table).
$ws.Rows.Item(6).RowHeight = 119

# Update the selection to match where the author was last working.
$ws.Range("P6:S6").Select()
